# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3156
$ws1.Range("F3").Value = 545
$ws1.Range("F4").Value = 1109
$ws1.Range("F5").Value = 98
$ws1.Range("F6").Value = 52
$ws1.Range("F8").Value = 43
$ws1.Range("F9").Value = 1141
$ws1.Range("F10").Value = 15982
$ws1.Range("F11").Value = 257
$ws1.Range("F14").Value = 6253
$ws1.Range("F15").Value = 631
$ws1.Range("F23").Value = 24
$ws1.Range("F29").Value = 5021
$ws1.Range("F31").Value = 11179
$ws1.Range("F38").Value = 74

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3156
$ws4.Range("F3").Value = 545
$ws4.Range("F4").Value = 1109
$ws4.Range("F5").Value = 98
$ws4.Range("F6").Value = 52
$ws4.Range("F8").Value = 43
$ws4.Range("F9").Value = 1141
$ws4.Range("F10").Value = 15982
$ws4.Range("F11").Value = 257
$ws4.Range("F14").Value = 6253
$ws4.Range("F15").Value = 631
$ws4.Range("F23").Value = 24
$ws4.Range("F29").Value = 5021
$ws4.Range("F32").Value = 11179
$ws4.Range("F39").Value = 74

$wb.Save()
